# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Update DAMSLTag (col I) / DialogAct (col J)
# values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 23;  I = "%";  J = "Uninterpretable" },
    @{ Row = 40;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 46;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 52;  I = "%";  J = "Uninterpretable" },
    @{ Row = 56;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 63;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 64;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 78;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 80;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 101; I = "aa"; J = "Agree/Accept" },
    @{ Row = 117; I = "%";  J = "Uninterpretable" },
    @{ Row = 123; I = "%";  J = "Uninterpretable" },
    @{ Row = 134; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 154; I = "ba"; J = "Appreciation" },
    @{ Row = 163; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 173; I = "sv"; J = "Statement-opinion" },
    @{ Row = 176; I = "sv"; J = "Statement-opinion" },
    @{ Row = 191; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 212; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 219; I = "sv"; J = "Statement-opinion" },
    @{ Row = 248; I = "sv"; J = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
